$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix "kcal/cap" -> "kcal/capita/day" in header W2
$ws.Range("W2").Value = "Total food supply (kcal/capita/day), 2019"

# Fix mojibake accented characters in the Regional Economic Communities note (A103)
$ws.Range("A103").Value = 'Regional Economic Communities:CEN-SAD = "Community of Sahel-Saharan States";COMESA = "Common Market for Eastern and Southern Africa";EAC = "East African Community";ECCAS = "Economic Community of Central African States";ECOWAS = "Economic Community of West African States";IGAD = "Intergovernmental Authority on Development";SADC = "Southern African Development Community";UMA = "Arab Maghreb Union";PALOP = "Países Africanos de Língua Oficial Portuguesa";ASEAN = "Association of Southeast Asian Nations";MERCOSUR = "Mercado Común del Sur".EU27 = "European Union (27 members)".OECD = "Organisation for Economic Co-operation and Development".'

# --- Header row 2: shift the column-group boundary from V to W ---
# (V2 becomes a "continuation" style like U2; W2 and X2 become "group start" style like old V2)
$ws.Range("V2").Copy()
$ws.Range("X2").PasteSpecial(-4122)
$ws.Range("V2").Copy()
$ws.Range("W2").PasteSpecial(-4122)
$ws.Range("U2").Copy()
$ws.Range("V2").PasteSpecial(-4122)

# --- Data rows 3-98: same column-group boundary shift, row by row ---
for ($row = 3; $row -le 98; $row++) {
    $ws.Range("V$row").Copy()
    $ws.Range("X$row").PasteSpecial(-4122)
    $ws.Range("V$row").Copy()
    $ws.Range("W$row").PasteSpecial(-4122)
    $ws.Range("U$row").Copy()
    $ws.Range("V$row").PasteSpecial(-4122)
}

# --- Update aggregate values in rows 97 and 98 (re-aggregated totals) ---
$row97 = @{
    "C" = 1469681.5544
    "D" = 1447357.2265
    "E" = 1960.5377
    "F" = 589939.822
    "G" = 176730.3
    "H" = 152891.6
    "I" = 336207.31
    "J" = 17588.0624
    "K" = 450626.9787
    "L" = 683239.7147
    "M" = 419208.5422
    "N" = 73123.8187
    "O" = 62472.7966
    "P" = 8916.6135
    "Q" = 10864.0401
    "R" = 141345.74
    "S" = 134349.396
    "T" = 117970.283
    "U" = 21450.389
    "V" = 5076.28
    "W" = 687.232002851554
    "X" = 498330.754
    "Y" = 23482.256
    "AC" = 1955.63258
    "AD" = 826.12345
    "AE" = 464.53843
}
foreach ($col in $row97.Keys) {
    $ws.Range($col + "97").Value = $row97[$col]
}

$row98 = @{
    "C" = 630368.0
    "D" = 602983.0
    "E" = 1715.6669
    "F" = 191387.4
    "G" = 92343.4
    "H" = 82169.9
    "I" = 178975.14
    "J" = 15057.7844
    "K" = 257334.9794
    "L" = 652647.1246
    "M" = 229450.6068
    "N" = 38978.0909
    "O" = 11846.3465
    "P" = 33918.9957
    "Q" = 29131.8842
    "R" = 183786.719
    "S" = 179113.203
    "T" = 139259.415
    "U" = 41898.345
    "V" = 2632.5
    "W" = 491.214832700331
    "X" = 196069.809
    "Y" = 6248.307
    "Z" = 6499.31276
    "AA" = 744.69
    "AC" = 7213.03743
    "AD" = 2525.48674
    "AE" = 976.24906
}
foreach ($col in $row98.Keys) {
    $ws.Range($col + "98").Value = $row98[$col]
}
